$wb = $excel.ActiveWorkbook

# 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 919
$ws1.Range("F6").Value = 2219

# 演出 sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14

# 全部类型 sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 14
$ws4.Range("F7").Value = 919
$ws4.Range("F8").Value = 2220
